# Auto-generated edit script: updates profit-calculation cells (columns H-N)
# on the Leve profit sheets, refreshing them with newer market-price snapshots
# (per the scheduled market-data refresh run).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3431.0344
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996
$ws.Range("H67").Value = 3431.0344
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216
$ws.Range("H76").Value = 33336216
$ws.Range("I76").Value = 45457070
$ws.Range("J76").Value = 3862.5
$ws.Range("K76").Value = 45457070
$ws.Range("L76").Value = 3862.5
$ws.Range("M76").Value = -45456755
$ws.Range("N76").Value = -4492.5
$ws.Range("H79").Value = 33336216
$ws.Range("I79").Value = 45457070
$ws.Range("J79").Value = 3862.5
$ws.Range("K79").Value = 45457070
$ws.Range("L79").Value = 3862.5
$ws.Range("M79").Value = -45455978
$ws.Range("N79").Value = -6046.5
$ws.Range("H111").Value = 610.5
$ws.Range("I111").Value = 511.66666
$ws.Range("K111").Value = 1534.99998
$ws.Range("M111").Value = 1532.00002
$ws.Range("H132").Value = 2472.5217
$ws.Range("I132").Value = 1755.8096
$ws.Range("K132").Value = 5267.4288
$ws.Range("M132").Value = -2737.4288
$ws.Range("H137").Value = 441495.4
$ws.Range("J137").Value = 807640.75
$ws.Range("L137").Value = 2422922.25
$ws.Range("N137").Value = -2428022.25
$ws.Range("H138").Value = 1918.1578
$ws.Range("I138").Value = 1511.8182
$ws.Range("J138").Value = 2476.875
$ws.Range("K138").Value = 4535.4546
$ws.Range("L138").Value = 7430.625
$ws.Range("M138").Value = 604.5454
$ws.Range("N138").Value = -17710.625
$ws.Range("H141").Value = 4375.5
$ws.Range("I141").Value = 3661.3
$ws.Range("J141").Value = 6161
$ws.Range("K141").Value = 10983.9
$ws.Range("L141").Value = 18483
$ws.Range("M141").Value = -5803.900000000001
$ws.Range("N141").Value = -28843

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6187.129
$ws.Range("I32").Value = 3217.2273
$ws.Range("J32").Value = 13446.889
$ws.Range("K32").Value = 3217.2273
$ws.Range("L32").Value = 13446.889
$ws.Range("M32").Value = -2930.2273
$ws.Range("N32").Value = -14020.889
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 235563
$ws.Range("I20").Value = 261575.67
$ws.Range("K20").Value = 261575.67
$ws.Range("M20").Value = -261328.67
$ws.Range("H94").Value = 2452.3333
$ws.Range("I94").Value = 2490.88
$ws.Range("J94").Value = 2259.6
$ws.Range("K94").Value = 2490.88
$ws.Range("L94").Value = 2259.6
$ws.Range("M94").Value = -2039.88
$ws.Range("N94").Value = -3161.6
$ws.Range("H105").Value = 64779.312
$ws.Range("I105").Value = 126281.75
$ws.Range("J105").Value = 3276.875
$ws.Range("K105").Value = 126281.75
$ws.Range("L105").Value = 3276.875
$ws.Range("M105").Value = -124534.75
$ws.Range("N105").Value = -6770.875
$ws.Range("H132").Value = 29490.877
$ws.Range("J132").Value = 29490.877
$ws.Range("L132").Value = 29490.877
$ws.Range("N132").Value = -39610.877
$ws.Range("H134").Value = 3934.484
$ws.Range("I134").Value = 3165.5186
$ws.Range("J134").Value = 9125
$ws.Range("K134").Value = 9496.5558
$ws.Range("L134").Value = 27375
$ws.Range("M134").Value = -6961.5558
$ws.Range("N134").Value = -32445
$ws.Range("H135").Value = 105163.336
$ws.Range("J135").Value = 105163.336
$ws.Range("L135").Value = 105163.336
$ws.Range("N135").Value = -115303.336
$ws.Range("H137").Value = 86666.664
$ws.Range("J137").Value = 86666.664
$ws.Range("L137").Value = 86666.664
$ws.Range("N137").Value = -96866.664
$ws.Range("H138").Value = 99760.75
$ws.Range("J138").Value = 99760.75
$ws.Range("L138").Value = 99760.75
$ws.Range("N138").Value = -110040.75
$ws.Range("H140").Value = 43500
$ws.Range("J140").Value = 43500
$ws.Range("L140").Value = 43500
$ws.Range("N140").Value = -53860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 804.2222
$ws.Range("I107").Value = 804.2222
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 804.2222
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1115.7778
$ws.Range("N107").Value = $null
$ws.Range("H122").Value = 2449.6155
$ws.Range("I122").Value = 2133
$ws.Range("J122").Value = 3162
$ws.Range("K122").Value = 6399
$ws.Range("L122").Value = 9486
$ws.Range("M122").Value = -3949
$ws.Range("N122").Value = -14386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7427634.5
$ws.Range("I4").Value = 6250464
$ws.Range("J4").Value = 16845000
$ws.Range("K4").Value = 18751392
$ws.Range("L4").Value = 50535000
$ws.Range("M4").Value = -18751280
$ws.Range("N4").Value = -50535224
$ws.Range("H18").Value = 12049.667
$ws.Range("I18").Value = 13305.875
$ws.Range("K18").Value = 39917.625
$ws.Range("M18").Value = -39748.625
$ws.Range("H48").Value = 999
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 999
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 2997
$ws.Range("M48").Value = $null
$ws.Range("N48").Value = -3497
$ws.Range("H50").Value = 2038.8572
$ws.Range("I50").Value = 1318
$ws.Range("J50").Value = 3000
$ws.Range("K50").Value = 3954
$ws.Range("L50").Value = 9000
$ws.Range("M50").Value = -3473
$ws.Range("N50").Value = -9962
$ws.Range("H53").Value = 2038.8572
$ws.Range("I53").Value = 1318
$ws.Range("J53").Value = 3000
$ws.Range("K53").Value = 3954
$ws.Range("L53").Value = 9000
$ws.Range("M53").Value = -3473
$ws.Range("N53").Value = -9962
$ws.Range("H70").Value = 1931.6666
$ws.Range("I70").Value = 1931.6666
$ws.Range("K70").Value = 5794.9998
$ws.Range("M70").Value = -5479.9998
$ws.Range("H73").Value = 1931.6666
$ws.Range("I73").Value = 1931.6666
$ws.Range("K73").Value = 5794.9998
$ws.Range("M73").Value = -4702.9998
$ws.Range("H75").Value = 231.25
$ws.Range("I75").Value = 217.5
$ws.Range("K75").Value = 652.5
$ws.Range("M75").Value = 345.5
$ws.Range("H78").Value = 231.25
$ws.Range("I78").Value = 217.5
$ws.Range("K78").Value = 1957.5
$ws.Range("M78").Value = 3034.5
$ws.Range("H92").Value = 243.70589
$ws.Range("I92").Value = 323.8
$ws.Range("K92").Value = 971.4000000000001
$ws.Range("M92").Value = 276.5999999999999
$ws.Range("H132").Value = 7587
$ws.Range("I132").Value = 1974.5
$ws.Range("J132").Value = 8834.223
$ws.Range("K132").Value = 17770.5
$ws.Range("L132").Value = 79508.007
$ws.Range("M132").Value = -15240.5
$ws.Range("N132").Value = -84568.007
$ws.Range("H133").Value = 7000
$ws.Range("I133").Value = 6000
$ws.Range("K133").Value = 18000
$ws.Range("M133").Value = -12940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 24069.334
$ws.Range("J52").Value = 23994.4
$ws.Range("L52").Value = 23994.4
$ws.Range("N52").Value = -24512.4
$ws.Range("H58").Value = 26498.75
$ws.Range("J58").Value = 26498.75
$ws.Range("L58").Value = 26498.75
$ws.Range("N58").Value = -27052.75
$ws.Range("H105").Value = 31375
$ws.Range("J105").Value = 31375
$ws.Range("L105").Value = 31375
$ws.Range("N105").Value = -38363
$ws.Range("H122").Value = 683336.2
$ws.Range("I122").Value = 876955.4
$ws.Range("J122").Value = 5669
$ws.Range("K122").Value = 2630866.2
$ws.Range("L122").Value = 17007
$ws.Range("M122").Value = -2628416.2
$ws.Range("N122").Value = -21907
$ws.Range("H132").Value = 3868.1538
$ws.Range("I132").Value = 3032.75
$ws.Range("J132").Value = 5204.8
$ws.Range("K132").Value = 9098.25
$ws.Range("L132").Value = 15614.4
$ws.Range("M132").Value = -6568.25
$ws.Range("N132").Value = -20674.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2214
$ws.Range("I82").Value = 2240.75
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 2240.75
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -1879.75
$ws.Range("N82").Value = -2722
$ws.Range("H85").Value = 2214
$ws.Range("I85").Value = 2240.75
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 2240.75
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -992.75
$ws.Range("N85").Value = -4496
$ws.Range("H106").Value = 14663.333
$ws.Range("J106").Value = 14663.333
$ws.Range("L106").Value = 14663.333
$ws.Range("N106").Value = -17187.333
$ws.Range("H132").Value = 1922.4
$ws.Range("I132").Value = 1770.25
$ws.Range("K132").Value = 5310.75
$ws.Range("M132").Value = -2780.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 94420.57000000001
$ws.Range("J46").Value = 94420.57000000001
$ws.Range("L46").Value = 94420.57000000001
$ws.Range("N46").Value = -94882.57000000001
$ws.Range("H81").Value = 2133.1667
$ws.Range("I81").Value = 2099.75
$ws.Range("J81").Value = 2200
$ws.Range("K81").Value = 4199.5
$ws.Range("L81").Value = 4400
$ws.Range("M81").Value = -3138.5
$ws.Range("N81").Value = -6522
$ws.Range("H84").Value = 2133.1667
$ws.Range("I84").Value = 2099.75
$ws.Range("J84").Value = 2200
$ws.Range("K84").Value = 20997.5
$ws.Range("L84").Value = 22000
$ws.Range("M84").Value = -15693.5
$ws.Range("N84").Value = -32608
$ws.Range("H113").Value = 1654.1052
$ws.Range("I113").Value = 1568.7778
$ws.Range("K113").Value = 4706.3334
$ws.Range("M113").Value = -2536.3334
$ws.Range("H134").Value = 94420.57000000001
$ws.Range("J134").Value = 94420.57000000001
$ws.Range("L134").Value = 283261.71
$ws.Range("N134").Value = -288331.71
$ws.Range("H136").Value = 782.2222
$ws.Range("I136").Value = 473.75
$ws.Range("K136").Value = 1421.25
$ws.Range("M136").Value = 1128.75
